$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 30")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# The run "QC_Aci_data.Rdata" (inside the "(QC_Aci_data.Rdata)" paragraph)
# starts at character 20 of the shape's text and is 17 characters long.
# It needs to become 4 separate runs: "QC", "_", "ACi", "_data.Rdata"
# (the only real content change is "Aci" -> "ACi").

$runQC = $tr.Characters(20, 2)
$runQC.Text = "QC"

$runUnderscore = $tr.Characters(22, 1)
$runUnderscore.Text = "_"

$runACi = $tr.Characters(23, 3)
$runACi.Text = "ACi"

$runRest = $tr.Characters(26, 11)
$runRest.Text = "_data.Rdata"
